$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=19; I='sv'; J='Statement-opinion'},
    @{Row=30; I='sv'; J='Statement-opinion'},
    @{Row=44; I='sv'; J='Statement-opinion'},
    @{Row=45; I='sd'; J='Statement-non-opinion'},
    @{Row=62; I='sd'; J='Statement-non-opinion'},
    @{Row=63; I='sv'; J='Statement-opinion'},
    @{Row=75; I='sv'; J='Statement-opinion'},
    @{Row=109; I='%'; J='Uninterpretable'},
    @{Row=111; I='sd'; J='Statement-non-opinion'},
    @{Row=124; I='sv'; J='Statement-opinion'},
    @{Row=127; I='ba'; J='Appreciation'},
    @{Row=151; I='sd'; J='Statement-non-opinion'},
    @{Row=166; I='sv'; J='Statement-opinion'},
    @{Row=175; I='aa'; J='Agree/Accept'},
    @{Row=179; I='sd'; J='Statement-non-opinion'},
    @{Row=183; I='sv'; J='Statement-opinion'},
    @{Row=205; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=207; I='sv'; J='Statement-opinion'},
    @{Row=222; I='sd'; J='Statement-non-opinion'},
    @{Row=242; I='sd'; J='Statement-non-opinion'},
    @{Row=252; I='sv'; J='Statement-opinion'},
    @{Row=255; I='aa'; J='Agree/Accept'},
    @{Row=257; I='sd'; J='Statement-non-opinion'},
    @{Row=260; I='sd'; J='Statement-non-opinion'},
    @{Row=265; I='sv'; J='Statement-opinion'},
    @{Row=271; I='sv'; J='Statement-opinion'},
    @{Row=274; I='aa'; J='Agree/Accept'},
    @{Row=294; I='sd'; J='Statement-non-opinion'},
    @{Row=304; I='sd'; J='Statement-non-opinion'},
    @{Row=339; I='aa'; J='Agree/Accept'},
    @{Row=350; I='sd'; J='Statement-non-opinion'},
    @{Row=354; I='sd'; J='Statement-non-opinion'},
    @{Row=358; I='sv'; J='Statement-opinion'},
    @{Row=364; I='aa'; J='Agree/Accept'},
    @{Row=366; I='sv'; J='Statement-opinion'},
    @{Row=367; I='sd'; J='Statement-non-opinion'},
    @{Row=372; I='%'; J='Uninterpretable'},
    @{Row=389; I='sv'; J='Statement-opinion'},
    @{Row=393; I='sv'; J='Statement-opinion'},
    @{Row=399; I='sd'; J='Statement-non-opinion'},
    @{Row=410; I='ba'; J='Appreciation'},
    @{Row=419; I='sv'; J='Statement-opinion'},
    @{Row=420; I='sd'; J='Statement-non-opinion'},
    @{Row=455; I='sd'; J='Statement-non-opinion'},
    @{Row=460; I='sd'; J='Statement-non-opinion'},
    @{Row=467; I='sd'; J='Statement-non-opinion'},
    @{Row=472; I='sv'; J='Statement-opinion'},
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
